$wb = $excel.ActiveWorkbook

# --- "Repayment schedule" sheet: insert a new (blank) column before column N ---
# to make room for a new "Late" sub-breakdown column used by the
# Loan RBI / Variable Instalments periodic schedule.
$ws = $wb.Worksheets.Item("Repayment schedule")

# Capture the width of the column immediately to the left (M) so the
# newly inserted column N inherits the same formatting/width Excel would
# naturally carry over on a column insert.
$leftWidth = $ws.Columns.Item(13).ColumnWidth

$ws.Columns.Item(14).Insert()
$ws.Columns.Item(14).ColumnWidth = $leftWidth

# Update the active selection/view on the sheet to reflect where the
# user ended up working after the edit.
$ws.Range("R7").Select()

# Make "Repayment schedule" the active tab of the workbook.
$ws.Activate()
